$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cryptocurrency price/volume data per latest scrape
$ws.Range('D2').Value = '63.739.50'
$ws.Range('E2').Value = '  -4.90%  '

$ws.Range('D3').Value = '3.498.31'
$ws.Range('E3').Value = '  -2.02%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.22%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '391.66'
$ws.Range('E5').Value = '  -5.90%  '

$ws.Range('E6').Value = '  -5.45%  '

$ws.Range('D7').Value = '3.483.80'
$ws.Range('E7').Value = '  -2.09%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  -9.85%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.669'
$ws.Range('E10').Value = '  -12.87%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.147'
$ws.Range('E11').Value = '  -15.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000341'
$ws.Range('E12').Value = '  +3.33%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.71'
$ws.Range('E13').Value = '  -8.36%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.22'
$ws.Range('E14').Value = '  -6.76%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.998.00'
$ws.Range('E15').Value = '  -3.47%  '

$ws.Range('E16').Value = '  -3.21%  '

$ws.Range('D17').Value = '3.502.12'
$ws.Range('E17').Value = '  -2.91%  '

$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.64'
$ws.Range('E18').Value = '  -8.65%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.41'
$ws.Range('E19').Value = '  +0.79%  '

$ws.Range('D20').Value = '63.623.23'
$ws.Range('E20').Value = '  -5.04%  '

$ws.Range('E21').Value = '  -10.76%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '392.88'
$ws.Range('E22').Value = '  -13.83%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.73'
$ws.Range('E23').Value = '  +2.98%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.27'
$ws.Range('E24').Value = '  -9.11%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').Value = '  -9.04%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.23'
$ws.Range('E26').Value = '  +7.23%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.06'
$ws.Range('E27').Value = '  -4.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.97'
$ws.Range('E28').Value = '  -11.60%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.77'
$ws.Range('E29').Value = '  -13.40%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '11.86'
$ws.Range('E30').Value = '  -3.87%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.62'
$ws.Range('E31').Value = '  -5.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.109'
$ws.Range('E32').Value = '  -6.95%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.76'
$ws.Range('E33').Value = '  -8.19%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.153'
$ws.Range('E34').Value = '  -5.32%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '36.79'
$ws.Range('E36').Value = '  -10.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.85'
$ws.Range('E37').Value = '  -4.89%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0436'
$ws.Range('E38').Value = '  -11.52%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.990'
$ws.Range('E39').Value = '  -0.91%  '

$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.86'
$ws.Range('E40').Value = '  +24.47%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.69'
$ws.Range('E41').Value = '  +14.64%  '

$ws.Range('B42').Value = 'PEPE'
$ws.Range('C42').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D42').Value = '0.0₃0635'
$ws.Range('E42').Value = '  -12.39%  '

$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.131'
$ws.Range('E43').Value = '  -10.78%  '

$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '139.23'
$ws.Range('E44').Value = '  -6.53%  '

$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.97'
$ws.Range('E45').Value = '  +10.90%  '

$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.95'
$ws.Range('E46').Value = '  -1.04%  '

$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.06'
$ws.Range('E47').Value = '  -6.13%  '

$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.70'
$ws.Range('E48').Value = '  -11.44%  '

$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.47'
$ws.Range('E49').Value = '  -9.48%  '

$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.04'
$ws.Range('E50').Value = '  -6.30%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.275'
$ws.Range('E51').Value = '  -11.91%  '
